$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F26").Value = 62
$ws.Range("G26").Value = 2858.2
$ws.Range("F38").Value = 6
$ws.Range("G38").Value = 184.44
$ws.Range("F40").Value = 8
$ws.Range("G40").Value = 368.8
$ws.Range("F41").Value = 10
$ws.Range("G41").Value = 307.4
$ws.Range("B46").Value = 23748.92
$ws.Range("F48").Value = 12
$ws.Range("G48").Value = 2361.24
$ws.Range("F55").Value = 185
$ws.Range("G55").Value = 35684.65
$ws.Range("F58").Value = 33
$ws.Range("G58").Value = 1165.56
$ws.Range("F64").Value = 99
$ws.Range("G64").Value = 3864.96
$ws.Range("F83").Value = 64
$ws.Range("G83").Value = 2141.44
$ws.Range("B85").Value = 138179.05
$ws.Range("F158").Value = 46
$ws.Range("G158").Value = 1951.78
$ws.Range("F164").Value = 77
$ws.Range("G164").Value = 1484.56
$ws.Range("B165").Value = 61218.84
$ws.Range("F193").Value = 20
$ws.Range("G193").Value = 1857.8
$ws.Range("F198").Value = 73
$ws.Range("G198").Value = 8408.139999999999
$ws.Range("B199").Value = 53925
$ws.Range("B200").Value = 57756
$ws.Range("B204").Value = 36186.71
$ws.Range("F210").Value = 4
$ws.Range("G210").Value = 2885.04
$ws.Range("B211").Value = 21900.5
$ws.Range("F226").Value = 67
$ws.Range("G226").Value = 5026.34
$ws.Range("B234").Value = 32317.92
$ws.Range("F258").Value = 20
$ws.Range("G258").Value = 1033
$ws.Range("F271").Value = 32
$ws.Range("G271").Value = 1674.24
$ws.Range("B273").Value = 16374.75
$ws.Range("B298").Value = 41864
$ws.Range("C298").Value = 'HAM-THERMOSTEEL 1000 ML WITH PLAIN LID'
$ws.Range("F298").Value = 0
$ws.Range("G298").Value = 0
$ws.Range("B299").Value = 56449
$ws.Range("C299").Value = 'HAM-Thermosteel 1000 Ml With Plain Lid'
$ws.Range("F299").Value = 24
$ws.Range("G299").Value = 16128.96
$ws.Range("F300").Value = 14
$ws.Range("G300").Value = 11321.52
$ws.Range("B304").Value = 96157.89999999999
$ws.Range("F321").Value = 2
$ws.Range("G321").Value = 357.82
$ws.Range("B355").Value = 130781.75
$ws.Range("F371").Value = 15
$ws.Range("G371").Value = 1977.45
$ws.Range("F387").Value = 123
$ws.Range("G387").Value = 2843.76
$ws.Range("F400").Value = 2
$ws.Range("G400").Value = 271.1
$ws.Range("F402").Value = 0
$ws.Range("G402").Value = 0
$ws.Range("F415").Value = 82
$ws.Range("G415").Value = 14049.06
$ws.Range("F420").Value = 10
$ws.Range("G420").Value = 1431
$ws.Range("F425").Value = 218
$ws.Range("G425").Value = 8977.24
$ws.Range("B429").Value = 92745.14
$ws.Range("F475").Value = 23
$ws.Range("G475").Value = 2138.77
$ws.Range("F486").Value = 34
$ws.Range("G486").Value = 4246.94
$ws.Range("B487").Value = 32658.1
$ws.Range("F500").Value = 21
$ws.Range("G500").Value = 1176.21
$ws.Range("B503").Value = 22905.67
$ws.Range("B524").Value = 47097
$ws.Range("D524").Value = 112.28
$ws.Range("E524").Value = 134.16
$ws.Range("F524").Value = 18
$ws.Range("G524").Value = 2021.04
$ws.Range("B525").Value = 58047
$ws.Range("D525").Value = 105.54
$ws.Range("E525").Value = 126.1
$ws.Range("F525").Value = 69
$ws.Range("G525").Value = 7282.26
$ws.Range("F533").Value = 144
$ws.Range("G533").Value = 5358.24
$ws.Range("F535").Value = 133
$ws.Range("G535").Value = 3576.37
$ws.Range("B538").Value = 139674.78
$ws.Range("F548").Value = 14
$ws.Range("G548").Value = 611.1
$ws.Range("B562").Value = 12663.64
$ws.Range("F567").Value = 737
$ws.Range("G567").Value = 9507.299999999999
$ws.Range("F569").Value = 158
$ws.Range("G569").Value = 2546.96
$ws.Range("F570").Value = 153
$ws.Range("G570").Value = 4322.25
$ws.Range("F571").Value = 122
$ws.Range("G571").Value = 4239.5
$ws.Range("F572").Value = 121
$ws.Range("G572").Value = 2324.41
$ws.Range("B573").Value = 44777.38
$ws.Range("F629").Value = 22
$ws.Range("G629").Value = 1063.92
$ws.Range("B631").Value = 45677.34
$ws.Range("F642").Value = 152
$ws.Range("G642").Value = 9229.440000000001
$ws.Range("B644").Value = 49151
$ws.Range("C644").Value = 'NES-MAGGI Atta Noodles Masala 290g'
$ws.Range("D644").Value = 78.09999999999999
$ws.Range("E644").Value = 88.58
$ws.Range("F644").Value = 1
$ws.Range("G644").Value = 78.09999999999999
$ws.Range("B645").Value = 55667
$ws.Range("C645").Value = 'NES-Maggi Atta Noodles Masala 290G'
$ws.Range("D645").Value = 85.76000000000001
$ws.Range("E645").Value = 97.25
$ws.Range("F645").Value = 13
$ws.Range("G645").Value = 1114.88
$ws.Range("F649").Value = 165
$ws.Range("G649").Value = 4042.5
$ws.Range("B659").Value = 125924.72
$ws.Range("F684").Value = 56
$ws.Range("G684").Value = 4412.8
$ws.Range("F685").Value = 2
$ws.Range("G685").Value = 151.12
$ws.Range("B688").Value = 20126.06
$ws.Range("F690").Value = 29
$ws.Range("G690").Value = 2297.67
$ws.Range("F694").Value = 121
$ws.Range("G694").Value = 3155.68
$ws.Range("B709").Value = 76537.03
$ws.Range("F752").Value = 22
$ws.Range("G752").Value = 2872.1
$ws.Range("F754").Value = 14
$ws.Range("G754").Value = 380.8
$ws.Range("F755").Value = 55
$ws.Range("G755").Value = 1496
$ws.Range("B757").Value = 15508.83
$ws.Range("F847").Value = 64
$ws.Range("G847").Value = 5829.12
$ws.Range("F848").Value = 93
$ws.Range("G848").Value = 8156.1
$ws.Range("F851").Value = 66
$ws.Range("G851").Value = 9261.780000000001
$ws.Range("B852").Value = 31845.73
$ws.Range("F858").Value = 84
$ws.Range("G858").Value = 11180.4
$ws.Range("F861").Value = 53
$ws.Range("G861").Value = 11442.17
$ws.Range("F867").Value = 215
$ws.Range("G867").Value = 32492.95
$ws.Range("B874").Value = 171087.34
$ws.Range("F881").Value = 3
$ws.Range("G881").Value = 75.98999999999999
$ws.Range("F893").Value = 466
$ws.Range("G893").Value = 14035.92
$ws.Range("F894").Value = 341
$ws.Range("G894").Value = 26799.19
$ws.Range("F900").Value = 124
$ws.Range("G900").Value = 3944.44
$ws.Range("F902").Value = 64
$ws.Range("G902").Value = 3192.96
$ws.Range("B904").Value = 177469.93
$ws.Range("F932").Value = 1
$ws.Range("G932").Value = 949.9400000000001
$ws.Range("B941").Value = 22698.1
$ws.Range("B964").Value = 3813033.36
$ws.Range("B965").Value = 3813033.36
